# Rename blank-node identifiers in column A (rows 2-22) so that the
# old GUID-based prefix "a040bad41d1b44d5b512af38ea3674f3b" is replaced
# with the new prefix "2e7a724ca8dc4fedaeaeed2f6551c45bb", keeping the
# trailing numeric suffix intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldPrefix = "na040bad41d1b44d5b512af38ea3674f3b"
$newPrefix = "n2e7a724ca8dc4fedaeaeed2f6551c45bb"

for ($row = 2; $row -le 22; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $current = [string]$cell.Value2
    if ($current -ne $null -and $current.StartsWith($oldPrefix)) {
        $suffix = $current.Substring($oldPrefix.Length)
        $cell.Value2 = $newPrefix + $suffix
    }
}
